# Education_treatment_final.xlsx update
#
# The MSfE (Manufacturing Systems for Engineers) rows were actually meant to
# describe the "Microsoft ... Bridging Program" certificate track, so the
# certificate_name column (D) for those six rows is corrected, and those
# cells get a new explicit-black-font style (matching the rest of the sheet
# visually while using a dedicated style entry). Two new data rows are also
# appended for a new "FSCeF" (Financial Services Connections En Français)
# treatment arm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the certificate_name text for the MSfE rows (column D), rows 44-49 ---
$ws.Range("D44").Value = "Microsoft Skills Bridging Program"
$ws.Range("D45").Value = "Microsoft  Bridging Program"
$ws.Range("D46").Value = "Microsoft Skills Bridging Program"
$ws.Range("D47").Value = "Microsoft  Bridging Program"
$ws.Range("D48").Value = "Microsoft Skills Bridging Program"
$ws.Range("D49").Value = "Microsoft  Bridging Program"

# These six cells pick up a fresh font (explicit black RGB) rather than the
# default/inherited color, giving them a distinct cell style.
$ws.Range("D44:D49").Font.Color = 0

# --- Append the two new "FSCeF" treatment rows at the bottom of the table ---
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = "FSCeF"
$ws.Range("C62").Value = "Collège Boréal – ACCES Employment"
$ws.Range("D62").Value = "Financial Services Connections En Français"
$ws.Range("E62").Value = 1

$ws.Range("A63").Value = 62
$ws.Range("B63").Value = "FSCeF"
$ws.Range("C63").Value = "Collège Boréal (In collaboration with ACCES Employment)"
$ws.Range("D63").Value = "Financial Services Connections Certificate (En Français)"
$ws.Range("E63").Value = 2
